$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.593.52'
$ws.Range("E2").Value = '  +2.12%  '
$ws.Range("D3").Value = '1.688.84'
$ws.Range("E3").Value = '  +3.38%  '
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '220.87'
$ws.Range("E5").Value = '  +2.77%  '
$ws.Range("E6").Value = '  +0.42%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = '30.92'
$ws.Range("E8").Value = '  +4.43%  '
$ws.Range("E9").Value = '  +2.33%  '
$ws.Range("E10").Value = '  +2.04%  '
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("D12").Value = '1.931.91'
$ws.Range("E12").Value = '  +3.48%  '
$ws.Range("D13").Value = '10.85'
$ws.Range("E13").Value = '  +14.13%  '
$ws.Range("D14").Value = '0.620'
$ws.Range("E14").Value = '  +8.31%  '
$ws.Range("D15").Value = '1.688.27'
$ws.Range("E15").Value = '  +3.51%  '
$ws.Range("E16").Value = '  +3.03%  '
$ws.Range("D17").Value = '30.583.35'
$ws.Range("E17").Value = '  +2.06%  '
$ws.Range("D18").Value = '65.94'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").Value = '249.56'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("E22").Value = '  +3.23%  '
$ws.Range("D23").Value = '10.21'
$ws.Range("E23").Value = '  +5.79%  '
$ws.Range("D24").Value = '2.21'
$ws.Range("E24").Value = '  +3.95%  '
$ws.Range("D25").Value = '157.57'
$ws.Range("E25").Value = '  -1.52%  '
$ws.Range("E26").Value = '  +1.47%  '
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("D28").Value = '6.78'
$ws.Range("E28").Value = '  +2.29%  '
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '1.14'
$ws.Range("E31").Value = '  +1.22%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '3.50'
$ws.Range("E32").Value = '  +3.45%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '3.31'
$ws.Range("E33").Value = '  +2.71%  '
$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '1.512.06'
$ws.Range("E34").Value = '  +5.55%  '
$ws.Range("E35").Value = '  +5.11%  '
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("E37").Value = '  +4.36%  '
$ws.Range("B38").Value = 'Aave'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D38").Value = '79.67'
$ws.Range("E38").Value = '  +8.30%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.72'
$ws.Range("E39").Value = '  -5.02%  '
$ws.Range("D40").Value = '0.587'
$ws.Range("E40").Value = '  +5.38%  '
$ws.Range("E41").Value = '  +1.37%  '
$ws.Range("E42").Value = '  +2.62%  '
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("D44").Value = '0.0504'
$ws.Range("E44").Value = '  +1.78%  '
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("D46").Value = '0.997'
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("D47").Value = '52.59'
$ws.Range("E47").Value = '  -4.12%  '
$ws.Range("D48").Value = '1.824.43'
$ws.Range("E48").Value = '  +2.82%  '
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("D50").Value = '95.73'
$ws.Range("E50").Value = '  +6.09%  '
$ws.Range("E51").Value = '  +7.04%  '
